$p = $ppt.ActivePresentation
$hm = $p.HandoutMaster
$hm.Background.Fill.ForeColor.RGB = 16777215
